$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    @(100, 0.04000000000000001, 6.683801412582397, 6.015625),
    @(100, 0.36, 107.3484826087952, 105.28125),
    @(100, 0.04000000000000001, 7.325151681900024, 6.625),
    @(100, 0.04000000000000001, 6.466713666915894, 5.453125),
    @(6.93, 0.04000000000000001, 5.264625072479248, 3.765625),
    @(12.09, 0.36, 36.54512047767639, 35.78125),
    @(8.49, 1, 85.68369960784912, 76.90625),
    @(11.92, 1, 55.40087199211121, 52.15625),
    @(16.63, 1, 47.94355058670044, 44.859375),
    @(37.37, 0.36, 22.59253191947937, 21.375),
    @(39.78, 0.16, 10.75686597824097, 10.21875),
    @(37.37, 0.36, 21.32076287269592, 20.75),
    @(37.37, 0.36, 21.63731241226196, 21.046875),
    @(17.84, 0.36, 17.41600775718689, 15.5625),
    @(17.84, 0.36, 18.63977336883545, 17.265625),
    @(17.84, 0.36, 15.21477603912354, 14.625),
    @(17.84, 0.36, 15.69570708274841, 14.671875),
    @(17.84, 0.36, 21.96741771697998, 20.890625),
    @(39.06, 0.04000000000000001, 2.837246417999268, 2.34375),
    @(37.37, 0.36, 25.06726121902466, 24.34375),
    @(39.78, 0.16, 10.740394115448, 18.109375),
    @(39.78, 0.16, 10.9668345451355, 17.703125),
    @(39.78, 0.16, 12.74186491966248, 12.203125)
)

$startRow = 10
for ($i = 0; $i -lt $newData.Count; $i++) {
    $rowValues = $newData[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowValues[$c]
    }
}
